{"js": "const replacements = [\n  [\"73-51=22\", \"92-61=31\"],\n  [\"18+74=92\", \"77-34=43\"],\n  [\"65-16=49\", \"45+25=70\"],\n  [\"35+13=48\", \"93-17=76\"],\n  [\"22-11=11\", \"29-12=17\"],\n  [\"88-81=7\", \"87-6=81\"],\n  [\"53+41=94\", \"28+6=34\"],\n  [\"73+21=94\", \"81-55=26\"],\n  [\"28+9=37\", \"7+47=54\"],\n  [\"97-88=9\", \"15+66=81\"],\n  [\"4-4=0\", \"86-26=60\"],\n  [\"67-52=15\", \"82+7=89\"],\n  [\"67-63=4\", \"37-0=37\"],\n  [\"47+36=83\", \"99-18=81\"],\n  [\"43+53=96\", \"66-48=18\"],\n  [\"77-38=39\", \"15-9=6\"],\n  [\"93-82=11\", \"79-18=61\"],\n  [\"34+63=97\", \"80-30=50\"],\n  [\"81-43=38\", \"72-1=71\"],\n  [\"52+26=78\", \"69+22=91\"],\n  [\"40+0=40\", \"31+55=86\"],\n  [\"95-42=53\", \"26+10=36\"],\n  [\"66-26=40\", \"10+45=55\"],\n  [\"71-7=64\", \"14+69=83\"],\n  [\"42+57=99\", \"54+8=62\"],\n  [\"24+73=97\", \"35+50=85\"],\n  [\"39+34=73\", \"10+42=52\"],\n  [\"70-12=58\", \"35+22=57\"],\n  [\"16+50=66\", \"42-28=14\"],\n  [\"81+8=89\", \"58-5=53\"],\n  [\"51+39=90\", \"69-27=42\"],\n  [\"74-4=70\", \"71-2=69\"],\n  [\"25+31=56\", \"45-33=12\"],\n  [\"1+19=20\", \"67+18=85\"],\n  [\"61-6=55\", \"0+40=40\"],\n  [\"63+16=79\", \"53-32=21\"],\n  [\"79-72=7\", \"34+27=61\"],\n  [\"56+42=98\", \"35-29=6\"],\n  [\"52-24=28\", \"23+32=55\"],\n  [\"97-60=37\", \"60+20=80\"],\n  [\"10+40=50\", \"91-3=88\"],\n  [\"55-7=48\", \"17+42=59\"],\n  [\"91-49=42\", \"7+51=58\"],\n  [\"19+44=63\", \"10+53=63\"],\n  [\"75-24=51\", \"82-78=4\"],\n  [\"15+10=25\", \"38+16=54\"],\n  [\"31+20=51\", \"4+83=87\"],\n  [\"52+13=65\", \"3+9=12\"],\n  [\"59-45=14\", \"18+40=58\"],\n  [\"76-59=17\", \"18+33=51\"],\n  [\"15-7=8\", \"13-11=2\"],\n  [\"14+32=46\", \"1+1=2\"],\n  [\"48-47=1\", \"95-30=65\"],\n  [\"45+15=60\", \"30+19=49\"],\n  [\"48-18=30\", \"0+59=59\"],\n  [\"32-5=27\", \"27+38=65\"],\n  [\"40+23=63\", \"74-25=49\"],\n  [\"1+22=23\", \"90-67=23\"],\n  [\"39-38=1\", \"22+35=57\"],\n  [\"96-3=93\", \"67-46=21\"],\n  [\"25-16=9\", \"3+1=4\"],\n  [\"22+34=56\", \"70-44=26\"],\n  [\"70+21=91\", \"28+44=72\"],\n  [\"31-25=6\", \"20+69=89\"],\n  [\"7+73=80\", \"80+11=91\"],\n  [\"19+18=37\", \"31+58=89\"],\n  [\"84+1=85\", \"55+33=88\"],\n  [\"42-7=35\", \"60+33=93\"],\n  [\"78-23=55\", \"55+21=76\"],\n  [\"91-0=91\", \"74-44=30\"],\n  [\"58+17=75\", \"41-29=12\"],\n  [\"29+3=32\", \"56+37=93\"],\n  [\"0+55=55\", \"21+36=57\"],\n  [\"43+56=99\", \"30+43=73\"],\n  [\"41+1=42\", \"17+26=43\"],\n  [\"10+20=30\", \"42+1=43\"],\n  [\"78-73=5\", \"80-13=67\"],\n  [\"32+26=58\", \"86-48=38\"],\n  [\"54-32=22\", \"84+13=97\"],\n  [\"71-27=44\", \"21+28=49\"],\n  [\"3+42=45\", \"17+9=26\"],\n  [\"8-0=8\", \"55+43=98\"],\n  [\"8+1=9\", \"59+11=70\"],\n  [\"80-54=26\", \"4+53=57\"],\n  [\"64-34=30\", \"47-47=0\"],\n  [\"14+55=69\", \"55+39=94\"],\n  [\"11+65=76\", \"43-5=38\"],\n  [\"0+30=30\", \"38+6=44\"],\n  [\"46+20=66\", \"87+11=98\"],\n  [\"41+15=56\", \"78+15=93\"],\n  [\"88-67=21\", \"19+4=23\"],\n  [\"49-35=14\", \"7+88=95\"],\n  [\"55-23=32\", \"57+31=88\"],\n  [\"28-10=18\", \"70-49=21\"],\n  [\"2+21=23\", \"39-27=12\"],\n  [\"63-24=39\", \"52-32=20\"],\n  [\"78-35=43\", \"41-14=27\"],\n  [\"21-15=6\", \"55-11=44\"],\n  [\"27+45=72\", \"1+74=75\"],\n  [\"82-53=29\", \"4+25=29\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,(\"73-51=22\", \"92-61=31\")\n    ,(\"18+74=92\", \"77-34=43\")\n    ,(\"65-16=49\", \"45+25=70\")\n    ,(\"35+13=48\", \"93-17=76\")\n    ,(\"22-11=11\", \"29-12=17\")\n    ,(\"88-81=7\", \"87-6=81\")\n    ,(\"53+41=94\", \"28+6=34\")\n    ,(\"73+21=94\", \"81-55=26\")\n    ,(\"28+9=37\", \"7+47=54\")\n    ,(\"97-88=9\", \"15+66=81\")\n    ,(\"4-4=0\", \"86-26=60\")\n    ,(\"67-52=15\", \"82+7=89\")\n    ,(\"67-63=4\", \"37-0=37\")\n    ,(\"47+36=83\", \"99-18=81\")\n    ,(\"43+53=96\", \"66-48=18\")\n    ,(\"77-38=39\", \"15-9=6\")\n    ,(\"93-82=11\", \"79-18=61\")\n    ,(\"34+63=97\", \"80-30=50\")\n    ,(\"81-43=38\", \"72-1=71\")\n    ,(\"52+26=78\", \"69+22=91\")\n    ,(\"40+0=40\", \"31+55=86\")\n    ,(\"95-42=53\", \"26+10=36\")\n    ,(\"66-26=40\", \"10+45=55\")\n    ,(\"71-7=64\", \"14+69=83\")\n    ,(\"42+57=99\", \"54+8=62\")\n    ,(\"24+73=97\", \"35+50=85\")\n    ,(\"39+34=73\", \"10+42=52\")\n    ,(\"70-12=58\", \"35+22=57\")\n    ,(\"16+50=66\", \"42-28=14\")\n    ,(\"81+8=89\", \"58-5=53\")\n    ,(\"51+39=90\", \"69-27=42\")\n    ,(\"74-4=70\", \"71-2=69\")\n    ,(\"25+31=56\", \"45-33=12\")\n    ,(\"1+19=20\", \"67+18=85\")\n    ,(\"61-6=55\", \"0+40=40\")\n    ,(\"63+16=79\", \"53-32=21\")\n    ,(\"79-72=7\", \"34+27=61\")\n    ,(\"56+42=98\", \"35-29=6\")\n    ,(\"52-24=28\", \"23+32=55\")\n    ,(\"97-60=37\", \"60+20=80\")\n    ,(\"10+40=50\", \"91-3=88\")\n    ,(\"55-7=48\", \"17+42=59\")\n    ,(\"91-49=42\", \"7+51=58\")\n    ,(\"19+44=63\", \"10+53=63\")\n    ,(\"75-24=51\", \"82-78=4\")\n    ,(\"15+10=25\", \"38+16=54\")\n    ,(\"31+20=51\", \"4+83=87\")\n    ,(\"52+13=65\", \"3+9=12\")\n    ,(\"59-45=14\", \"18+40=58\")\n    ,(\"76-59=17\", \"18+33=51\")\n    ,(\"15-7=8\", \"13-11=2\")\n    ,(\"14+32=46\", \"1+1=2\")\n    ,(\"48-47=1\", \"95-30=65\")\n    ,(\"45+15=60\", \"30+19=49\")\n    ,(\"48-18=30\", \"0+59=59\")\n    ,(\"32-5=27\", \"27+38=65\")\n    ,(\"40+23=63\", \"74-25=49\")\n    ,(\"1+22=23\", \"90-67=23\")\n    ,(\"39-38=1\", \"22+35=57\")\n    ,(\"96-3=93\", \"67-46=21\")\n    ,(\"25-16=9\", \"3+1=4\")\n    ,(\"22+34=56\", \"70-44=26\")\n    ,(\"70+21=91\", \"28+44=72\")\n    ,(\"31-25=6\", \"20+69=89\")\n    ,(\"7+73=80\", \"80+11=91\")\n    ,(\"19+18=37\", \"31+58=89\")\n    ,(\"84+1=85\", \"55+33=88\")\n    ,(\"42-7=35\", \"60+33=93\")\n    ,(\"78-23=55\", \"55+21=76\")\n    ,(\"91-0=91\", \"74-44=30\")\n    ,(\"58+17=75\", \"41-29=12\")\n    ,(\"29+3=32\", \"56+37=93\")\n    ,(\"0+55=55\", \"21+36=57\")\n    ,(\"43+56=99\", \"30+43=73\")\n    ,(\"41+1=42\", \"17+26=43\")\n    ,(\"10+20=30\", \"42+1=43\")\n    ,(\"78-73=5\", \"80-13=67\")\n    ,(\"32+26=58\", \"86-48=38\")\n    ,(\"54-32=22\", \"84+13=97\")\n    ,(\"71-27=44\", \"21+28=49\")\n    ,(\"3+42=45\", \"17+9=26\")\n    ,(\"8-0=8\", \"55+43=98\")\n    ,(\"8+1=9\", \"59+11=70\")\n    ,(\"80-54=26\", \"4+53=57\")\n    ,(\"64-34=30\", \"47-47=0\")\n    ,(\"14+55=69\", \"55+39=94\")\n    ,(\"11+65=76\", \"43-5=38\")\n    ,(\"0+30=30\", \"38+6=44\")\n    ,(\"46+20=66\", \"87+11=98\")\n    ,(\"41+15=56\", \"78+15=93\")\n    ,(\"88-67=21\", \"19+4=23\")\n    ,(\"49-35=14\", \"7+88=95\")\n    ,(\"55-23=32\", \"57+31=88\")\n    ,(\"28-10=18\", \"70-49=21\")\n    ,(\"2+21=23\", \"39-27=12\")\n    ,(\"63-24=39\", \"52-32=20\")\n    ,(\"78-35=43\", \"41-14=27\")\n    ,(\"21-15=6\", \"55-11=44\")\n    ,(\"27+45=72\", \"1+74=75\")\n    ,(\"82-53=29\", \"4+25=29\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $oldText,\n        $false,\n        $true,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}"}
